# Auto-generated edit script: Add data for 2025-07-17
# Updates cumulative 2025 (column L) violent-crime counts across
# Citywide Totals, By Neighborhood pivot, and per-neighborhood sheets,
# plus a few minor retroactive corrections to 2023 (J) / 2024 (K) columns.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("L2").Value = 3592
$ws.Range("L3").Value = 3765
$ws.Range("J4").Value = 1868
$ws.Range("L4").Value = 935
$ws.Range("K5").Value = 588
$ws.Range("L5").Value = 223
$ws.Range("L6").Value = 3280
$ws.Range("J7").Value = 29344
$ws.Range("K7").Value = 27563
$ws.Range("L7").Value = 11795

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("L4").Value = 45
$ws.Range("L5").Value = 46
$ws.Range("L8").Value = 760
$ws.Range("L11").Value = 197
$ws.Range("L15").Value = 85
$ws.Range("L17").Value = 22
$ws.Range("L18").Value = 91
$ws.Range("L19").Value = 337
$ws.Range("L20").Value = 297
$ws.Range("L23").Value = 129
$ws.Range("L27").Value = 108
$ws.Range("L29").Value = 637
$ws.Range("L31").Value = 114
$ws.Range("L33").Value = 557
$ws.Range("L36").Value = 160
$ws.Range("L37").Value = 420
$ws.Range("L42").Value = 379
$ws.Range("L44").Value = 86
$ws.Range("L48").Value = 166
$ws.Range("L51").Value = 146
$ws.Range("L52").Value = 242
$ws.Range("L53").Value = 133
$ws.Range("L60").Value = 70
$ws.Range("K63").Value = 164
$ws.Range("L67").Value = 420
$ws.Range("L76").Value = 174
$ws.Range("L78").Value = 147
$ws.Range("L79").Value = 308
$ws.Range("L83").Value = 270
$ws.Range("J85").Value = 1190
$ws.Range("L85").Value = 603
$ws.Range("L86").Value = 88
$ws.Range("L89").Value = 166
$ws.Range("L91").Value = 169
$ws.Range("L94").Value = 140
$ws.Range("L96").Value = 118
$ws.Range("L97").Value = 101
$ws.Range("L99").Value = 199
$ws.Range("J101").Value = 29344
$ws.Range("K101").Value = 27563
$ws.Range("L101").Value = 11795

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range("L6").Value = 31
$ws.Range("L7").Value = 118

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range("L2").Value = 76
$ws.Range("L6").Value = 48
$ws.Range("L7").Value = 197

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range("L3").Value = 48
$ws.Range("L7").Value = 166

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("L3").Value = 251
$ws.Range("J4").Value = 73
$ws.Range("J7").Value = 1190
$ws.Range("L7").Value = 603

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("L3").Value = 74
$ws.Range("L7").Value = 242

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range("L2").Value = 41
$ws.Range("L3").Value = 34
$ws.Range("L6").Value = 45
$ws.Range("L7").Value = 133

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("L3").Value = 254
$ws.Range("L6").Value = 210
$ws.Range("L7").Value = 760

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("L3").Value = 107
$ws.Range("L4").Value = 9
$ws.Range("L7").Value = 270

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("L4").Value = 31
$ws.Range("L5").Value = 10
$ws.Range("L6").Value = 182
$ws.Range("L7").Value = 557

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("L2").Value = 129
$ws.Range("L3").Value = 133
$ws.Range("L6").Value = 119
$ws.Range("L7").Value = 420

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range("L3").Value = 82
$ws.Range("L7").Value = 199

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range("L2").Value = 40
$ws.Range("L7").Value = 114

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("L2").Value = 121
$ws.Range("L3").Value = 158
$ws.Range("L7").Value = 420

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("L2").Value = 192
$ws.Range("L3").Value = 244
$ws.Range("L6").Value = 159
$ws.Range("L7").Value = 637

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range("L6").Value = 74
$ws.Range("L7").Value = 166

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("L3").Value = 103
$ws.Range("L6").Value = 99
$ws.Range("L7").Value = 337

$ws = $wb.Worksheets.Item('Irving Park')
$ws.Range("L2").Value = 35
$ws.Range("L4").Value = 7
$ws.Range("L7").Value = 86

$ws = $wb.Worksheets.Item('River North')
$ws.Range("L6").Value = 81
$ws.Range("L7").Value = 174

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("L2").Value = 114
$ws.Range("L3").Value = 119
$ws.Range("L6").Value = 108
$ws.Range("L7").Value = 379

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range("L3").Value = 42
$ws.Range("L7").Value = 147

$ws = $wb.Worksheets.Item('Douglas')
$ws.Range("L3").Value = 50
$ws.Range("L6").Value = 34
$ws.Range("L7").Value = 129

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range("L2").Value = 60
$ws.Range("L3").Value = 71
$ws.Range("L7").Value = 169

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("L2").Value = 105
$ws.Range("L6").Value = 61
$ws.Range("L7").Value = 308

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("L2").Value = 95
$ws.Range("L3").Value = 94
$ws.Range("L7").Value = 297

$ws = $wb.Worksheets.Item('Calumet Heights')
$ws.Range("L4").Value = 9
$ws.Range("L6").Value = 15
$ws.Range("L7").Value = 91

$ws = $wb.Worksheets.Item('Burnside')
$ws.Range("L6").Value = 7
$ws.Range("L7").Value = 22

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range("L2").Value = 60
$ws.Range("L7").Value = 160

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range("L2").Value = 33
$ws.Range("L3").Value = 32
$ws.Range("L7").Value = 140

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range("L3").Value = 24
$ws.Range("L7").Value = 85

$ws = $wb.Worksheets.Item('West Town')
$ws.Range("L2").Value = 23
$ws.Range("L3").Value = 21
$ws.Range("L7").Value = 101

$ws = $wb.Worksheets.Item('Armour Square')
$ws.Range("L6").Value = 21
$ws.Range("L7").Value = 46

$ws = $wb.Worksheets.Item('Edgewater')
$ws.Range("L3").Value = 35
$ws.Range("L7").Value = 108

$ws = $wb.Worksheets.Item('Streeterville')
$ws.Range("L4").Value = 47
$ws.Range("L7").Value = 88

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range("L2").Value = 45
$ws.Range("L4").Value = 21
$ws.Range("L7").Value = 146

$ws = $wb.Worksheets.Item('Morgan Park')
$ws.Range("L2").Value = 21
$ws.Range("L7").Value = 70

$ws = $wb.Worksheets.Item('Archer Heights')
$ws.Range("L2").Value = 17
$ws.Range("L7").Value = 45
